# Apply cryptos list update (prices/volumes refreshed; two rank swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '67.808.01'
Set-TextCell $ws.Range('E2') '  +0.59%  '

Set-TextCell $ws.Range('D3') '3.810.26'
Set-TextCell $ws.Range('E3') '  +1.05%  '

Set-TextCell $ws.Range('E4') '  -0.07%  '

Set-TextCell $ws.Range('D5') '598.01'
Set-TextCell $ws.Range('E5') '  +0.75%  '

Set-TextCell $ws.Range('D6') '167.64'
Set-TextCell $ws.Range('E6') '  +1.02%  '

Set-TextCell $ws.Range('E8') '  -0.02%  '

Set-TextCell $ws.Range('D9') '0.161'
Set-TextCell $ws.Range('E9') '  +0.96%  '

Set-TextCell $ws.Range('D10') '6.31'
Set-TextCell $ws.Range('E10') '  -1.52%  '

Set-TextCell $ws.Range('E11') '  -0.27%  '

Set-TextCell $ws.Range('E12') '  -0.98%  '

Set-TextCell $ws.Range('D13') '36.09'

Set-TextCell $ws.Range('D14') '4.444.25'
Set-TextCell $ws.Range('E14') '  +0.95%  '

Set-TextCell $ws.Range('D15') '3.815.07'
Set-TextCell $ws.Range('E15') '  +0.77%  '

Set-TextCell $ws.Range('D16') '18.63'
Set-TextCell $ws.Range('E16') '  +5.01%  '

Set-TextCell $ws.Range('D17') '67.791.13'
Set-TextCell $ws.Range('E17') '  +0.45%  '

Set-TextCell $ws.Range('E18') '  +2.46%  '

Set-TextCell $ws.Range('D20') '462.33'
Set-TextCell $ws.Range('E20') '  +0.92%  '

Set-TextCell $ws.Range('E21') '  -5.90%  '

Set-TextCell $ws.Range('D23') '0.0000155'
Set-TextCell $ws.Range('E23') '  +1.22%  '

Set-TextCell $ws.Range('D24') '83.64'
Set-TextCell $ws.Range('E24') '  +0.37%  '

Set-TextCell $ws.Range('E25') '  +2.73%  '

Set-TextCell $ws.Range('E26') '  -1.43%  '

Set-TextCell $ws.Range('D27') '10.04'
Set-TextCell $ws.Range('E27') '  +0.04%  '

Set-TextCell $ws.Range('D28') '1.01'
Set-TextCell $ws.Range('E28') '  +0.67%  '

Set-TextCell $ws.Range('D29') '3.955.99'
Set-TextCell $ws.Range('E29') '  +0.99%  '

Set-TextCell $ws.Range('D30') '2.80'
Set-TextCell $ws.Range('E30') '  +1.22%  '

Set-TextCell $ws.Range('D31') '2.25'
Set-TextCell $ws.Range('E31') '  +4.26%  '

Set-TextCell $ws.Range('E32') '  +1.15%  '

Set-TextCell $ws.Range('D33') '29.76'
Set-TextCell $ws.Range('E33') '  -0.01%  '

Set-TextCell $ws.Range('B34') 'Binance-PegBSC-USD'
Set-TextCell $ws.Range('C34') 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws.Range('D34') '1.00'
Set-TextCell $ws.Range('E34') '  +0.25%  '

Set-TextCell $ws.Range('B35') 'Aptos'
Set-TextCell $ws.Range('C35') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range('D35') '9.11'
Set-TextCell $ws.Range('E35') '  -0.64%  '

Set-TextCell $ws.Range('D36') '3.749.24'
Set-TextCell $ws.Range('E36') '  +0.74%  '

Set-TextCell $ws.Range('E38') '  +0.81%  '

Set-TextCell $ws.Range('E39') '  +0.65%  '

Set-TextCell $ws.Range('E40') '  +0.85%  '

Set-TextCell $ws.Range('D41') '5.80'
Set-TextCell $ws.Range('E41') '  +1.00%  '

Set-TextCell $ws.Range('D42') '0.999'
Set-TextCell $ws.Range('E42') '  -0.11%  '

Set-TextCell $ws.Range('D44') '48.17'
Set-TextCell $ws.Range('E44') '  +2.75%  '

Set-TextCell $ws.Range('D45') '43.77'
Set-TextCell $ws.Range('E45') '  -0.84%  '

Set-TextCell $ws.Range('D46') '0.300'
Set-TextCell $ws.Range('E46') '  +0.18%  '

Set-TextCell $ws.Range('D47') '149.13'
Set-TextCell $ws.Range('E47') '  +2.23%  '

Set-TextCell $ws.Range('B48') 'Cosmos'
Set-TextCell $ws.Range('C48') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws.Range('D48') '8.35'
Set-TextCell $ws.Range('E48') '  -0.13%  '

Set-TextCell $ws.Range('B49') 'Bittensor'
Set-TextCell $ws.Range('C49') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws.Range('D49') '400.12'
Set-TextCell $ws.Range('E49') '  +2.10%  '

Set-TextCell $ws.Range('E50') '  -2.75%  '

Set-TextCell $ws.Range('D51') '26.75'
Set-TextCell $ws.Range('E51') '  +5.89%  '
